# Fruta / hortaliza, semanal
# Re-sequences the weekly rows (2..10) of the "Granada" sheet: the data in
# columns D (Fecha) and K..T (Variedad .. Kg/unidad) is reshuffled across
# rows while A..J (market/product identifiers) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each row, keyed by destination row number.
$rows = @{
    2  = @{ D = 44280; K = "Sin especificar"; L = "Primera"; M = 15;  N = 360000; O = 360000; P = 360000; Q = '$/bins (450 kilos)';       R = "Provincia del Elquí";   S = 800;  T = 450 }
    3  = @{ D = 44307; K = "Sin especificar"; L = "Primera"; M = 150; N = 16000;  O = 18000;  P = 17000;  Q = '$/caja 15 kilos granel';    R = "Región de O'Higgins";   S = 1133; T = 15  }
    4  = @{ D = 44320; K = "Wonderfull";      L = "Primera"; M = 12;  N = 250000; O = 260000; P = 255000; Q = '$/bins (400 kilos)';       R = "Provincia de Limarí";   S = 638;  T = 400 }
    5  = @{ D = 44285; K = "Wonderfull";      L = "Primera"; M = 8;   N = 280000; O = 300000; P = 290000; Q = '$/bins (400 kilos)';       R = "Provincia del Elquí";   S = 725;  T = 400 }
    6  = @{ D = 44266; K = "Wonderfull";      L = "Segunda"; M = 120; N = 4800;   O = 4800;   P = 4800;   Q = '$/bandeja 4 kilos';        R = "Provincia del Elquí";   S = 1200; T = 4   }
    7  = @{ D = 44266; K = "Wonderfull";      L = "Tercera"; M = 80;  N = 4000;   O = 4000;   P = 4000;   Q = '$/bandeja 4 kilos';        R = "Provincia del Elquí";   S = 1000; T = 4   }
    8  = @{ D = 44334; K = "Wonderfull";      L = "Primera"; M = 16;  N = 240000; O = 250000; P = 245000; Q = '$/bins (450 kilos)';       R = "Provincia de Limarí";   S = 544;  T = 450 }
    9  = @{ D = 44312; K = "Wonderfull";      L = "Primera"; M = 24;  N = 220000; O = 240000; P = 230000; Q = '$/bins (400 kilos)';       R = "Región de O'Higgins";   S = 575;  T = 400 }
    10 = @{ D = 44312; K = "Wonderfull";      L = "Primera"; M = 34;  N = 240000; O = 240000; P = 240000; Q = '$/bins (450 kilos)';       R = "Provincia del Elquí";   S = 533;  T = 450 }
}

foreach ($r in $rows.Keys) {
    $v = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $v.D   # D: Fecha
    $ws.Cells.Item($r, 11).Value = $v.K   # K: Variedad
    $ws.Cells.Item($r, 12).Value = $v.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $v.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $v.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $v.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $v.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $v.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $v.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $v.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $v.T   # T: Kg / unidad
}
